$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 78 (shifts existing rows 78-100 down to 79-101)
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new weekly price record
$ws.Cells.Item(78, 1).Value = 3
$ws.Cells.Item(78, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(78, 3).Value = "Coquimbo"
$ws.Cells.Item(78, 4).Value = 44524
$ws.Cells.Item(78, 5).Value = 5
$ws.Cells.Item(78, 6).Value = 100112026
$ws.Cells.Item(78, 7).Value = "Haba"
$ws.Cells.Item(78, 8).Value = "Sin especificar"
$ws.Cells.Item(78, 9).Value = "Primera"
$ws.Cells.Item(78, 10).Value = 80
$ws.Cells.Item(78, 11).Value = 8000
$ws.Cells.Item(78, 12).Value = 8500
$ws.Cells.Item(78, 13).Value = 8250
$ws.Cells.Item(78, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(78, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(78, 16).Value = 330
$ws.Cells.Item(78, 17).Value = 25
$ws.Cells.Item(78, 18).Value = "Hortaliza"

# Match the date cell's number format to the rest of column D (dates)
$ws.Cells.Item(78, 4).NumberFormat = $ws.Cells.Item(79, 4).NumberFormat()
